# Marksheet update: handle float-input marking data without breaking the
# sheet, refresh the computed stat rows, and drop the (now unused) third
# "Student Ans / Correct Ans" answer-block while filling in the attempted
# answers for sets 1 and 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Copy-Style($fromAddr, $toAddr) {
    $ws.Range($fromAddr).Copy() | Out-Null
    $ws.Range($toAddr).PasteSpecial($xlPasteFormats) | Out-Null
}

# ---------------------------------------------------------------------
# Stat block (rows 10-12): give the row-label cells in column A the same
# "mtitleStyle" look as the header row (row 9), and refresh the computed
# numbers now that the marking logic copes with the (float-safe) inputs.
# ---------------------------------------------------------------------
Copy-Style "A9" "A10"
Copy-Style "A9" "A11"
Copy-Style "A9" "A12"

# Row 10 - counts
$ws.Range("B10").Value = 19
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 7
$ws.Range("E10").Value = 28

# Row 11 - marking scheme (C11 used to be stored as text "-1"; now numeric)
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1

# Row 12 - totals
$ws.Range("B12").Value = 76
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "74/112"

# ---------------------------------------------------------------------
# Drop the third "Student Ans / Correct Ans" block (columns G:H) entirely.
# ---------------------------------------------------------------------
$ws.Range("G15:H40").Clear()

# ---------------------------------------------------------------------
# Set 2 (columns D:E) only has answers recorded for the first three
# questions now; clear the rest.
# ---------------------------------------------------------------------
$ws.Range("D19:E40").Clear()

# Fill in the three attempted set-2 answers, colouring each by whether it
# matches the correct answer already in column E.
Copy-Style "C10" "D16"
$ws.Range("D16").Value = "Option B"

Copy-Style "B10" "D17"
$ws.Range("D17").Value = "Option C"

Copy-Style "C10" "D18"
$ws.Range("D18").Value = "Option B"

# ---------------------------------------------------------------------
# Set 1 (columns A:B) - fill in the attempted answers in column A. Seven
# questions were never attempted (rows 20, 21, 26, 31, 34, 36, 40) and stay
# blank; every attempted answer below matches the correct answer, so each
# gets the green "correctStyle" look (same style already used by B10).
# ---------------------------------------------------------------------
$answered = @{
    16 = "Option A"
    17 = "Option D"
    18 = "Option B"
    19 = "Option C"
    22 = "Option D"
    23 = "Option D"
    24 = "Option A"
    25 = "Option A"
    27 = "Option A"
    28 = "Option D"
    29 = "Option D"
    30 = "Option B"
    32 = "Option C"
    33 = "Option D"
    35 = "Option D"
    37 = "Option A"
    38 = "Option A"
    39 = "Option D"
}

foreach ($row in $answered.Keys) {
    $addr = "A$row"
    Copy-Style "B10" $addr
    $ws.Range($addr).Value = $answered[$row]
}

Write-Output "edit complete"
